$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.253.51'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.539.72'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.68'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.59'
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.539.91'
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("E10").Value = '  -4.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '8.05'
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.415'
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.138.91'
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("E14").Value = '  -3.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.48'
$ws.Range("E15").Value = '  -3.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.541.16'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.370.60'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.92'
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("E20").Value = '  -2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.03'
$ws.Range("E21").Value = '  -2.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '426.37'
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.53'
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.682.40'
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.35'
$ws.Range("E28").Value = '  -5.12%  '
$ws.Range("E29").Value = '  -4.25%  '
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.49'
$ws.Range("E32").Value = '  -6.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.157'
$ws.Range("E33").Value = '  -7.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.39'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.528.24'
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.87'
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("E39").Value = '  -4.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '169.79'
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0863'
$ws.Range("E42").Value = '  -3.27%  '
$ws.Range("E43").Value = '  -4.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.894'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("E45").Value = '  -9.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.47'
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("E47").Value = '  -7.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.12'
$ws.Range("E48").Value = '  -7.76%  '
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("E50").Value = '  -3.77%  '
$ws.Range("E51").Value = '  -3.96%  '
